# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price-column ("D") values are plain text in this sheet (not numbers), so for
# any replacement value that Excel would otherwise auto-parse as a number we
# force the cell to Text format first to preserve the literal digits/zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.206.13'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '2.303.22'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.75'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.00'
$ws.Range('E6').Value = '  +7.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.525'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  +6.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.82'
$ws.Range('E10').Value = '  +3.24%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0808'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.93'
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').Value = '2.660.17'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.06'
$ws.Range('E16').Value = '  +4.21%  '
$ws.Range('D17').Value = '2.302.32'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.801'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '43.163.26'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.94'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').Value = '0.0₃0921'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.16'
$ws.Range('E22').Value = '  +4.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.74'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '239.92'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  +3.36%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.59'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.67'
$ws.Range('E28').Value = '  +5.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.27'
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.55'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.10'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.37'
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.23'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.19'
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.54'
$ws.Range('E36').Value = '  +6.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0734'
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.50'
$ws.Range('E39').Value = '  +10.80%  '
$ws.Range('E40').Value = '  +3.81%  '
$ws.Range('E41').Value = '  +3.57%  '
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.53'
$ws.Range('E43').Value = '  +15.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0289'
$ws.Range('E44').Value = '  +2.81%  '
$ws.Range('D45').Value = '1.962.98'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.77'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.08'
$ws.Range('E47').Value = '  +6.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.17'
$ws.Range('E48').Value = '  +5.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '57.03'
$ws.Range('E49').Value = '  +6.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('E51').Value = '  +7.94%  '
